$d = $word.ActiveDocument
$d.Content.Find.Execute("list of countries", $true, $false, $false, $false, $false,
                         $true, 1, $false, "list of regions", 2)
